# Main manuscript tables: missing value formatting and rounding.
#
# Table 2 ("ATC | Name | N | AUC | Q value"):
#   - header cell "AUC" -> "round(AUC, 3)"
#   - AUC / Q value columns rounded to 3 decimal places
# Table 3 (positional_mapping / fine_mapping enrichment table):
#   - "NA" placeholders replaced with an em dash "—"

$d = $word.ActiveDocument

$emDash = [char]0x2014

# --- Table 2: AUC / Q value table -----------------------------------------
$tbl2 = $d.Tables.Item(2)

# Header: "AUC" -> "round(AUC, 3)"
$tbl2.Cell(1, 4).Range.Text = "round(AUC, 3)"

# Row 2: N05A / Antipsychotics
$tbl2.Cell(2, 4).Range.Text = "0.646"
$tbl2.Cell(2, 5).Range.Text = "0.010"

# Row 3: J01M / Quinolone Antibacterials
$tbl2.Cell(3, 4).Range.Text = "0.764"
$tbl2.Cell(3, 5).Range.Text = "0.010"

# Row 4: J01MA / Fluoroquinolones
$tbl2.Cell(4, 4).Range.Text = "0.780"
$tbl2.Cell(4, 5).Range.Text = "0.034"

# Row 5: N05 / Psycholeptics
$tbl2.Cell(5, 4).Range.Text = "0.587"
$tbl2.Cell(5, 5).Range.Text = "0.035"

# Row 6: S01A / Antiinfectives
$tbl2.Cell(6, 4).Range.Text = "0.633"
$tbl2.Cell(6, 5).Range.Text = "0.035"

# --- Table 3: positional_mapping / fine_mapping table ----------------------
$tbl3 = $d.Tables.Item(3)

# Row 2 (positional_mapping), column 2 (N06A (MR-MEGA)): "NA" -> em dash
$tbl3.Cell(2, 2).Range.Text = $emDash

# Row 3 (fine_mapping), columns 3-5 (N06A/N06AA/N06AB (EUR)): "NA" -> em dash
$tbl3.Cell(3, 3).Range.Text = $emDash
$tbl3.Cell(3, 4).Range.Text = $emDash
$tbl3.Cell(3, 5).Range.Text = $emDash
